# Add a new worksheet "Sheet3" at the end of the workbook, containing a short
# list of page titles used to validate social media links.
$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "Sheet3"

$ws.Range("A1").Value = "Sign In | LinkedIn"
$ws.Range("A2").Value = "Sauce Labs - Home"
$ws.Range("A3").Value = "Twitter"

$ws.Columns.Item(1).ColumnWidth = 40.5

# Match the page setup (paper size / orientation) used by the other sheets.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Mirror the author's saved selection/active cell on the new sheet.
[void]$ws.Range("A5").Select()
